# Weekly update: a new price-report row is inserted at row 109 (pushing the
# existing rows 109-212 down to 110-213, and appending the former row 212's
# data onto the newly created row 213), then the new row 109 is populated
# with this week's record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 109; everything below shifts
# down by one (row 109 -> 110, ..., row 212 -> 213).
$ws.Rows("109").Insert()

# Populate the newly inserted row 109 with the new weekly record.
$ws.Cells.Item(109, 1).Value = 5
$ws.Cells.Item(109, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(109, 3).Value = "Maule"
$ws.Cells.Item(109, 4).Value = 45271
$ws.Cells.Item(109, 5).Value = 7
$ws.Cells.Item(109, 6).Value = "Fruta"
$ws.Cells.Item(109, 7).Value = 100108
$ws.Cells.Item(109, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(109, 9).Value = 100108002
$ws.Cells.Item(109, 10).Value = "Mango"
$ws.Cells.Item(109, 11).Value = "Sin especificar"
$ws.Cells.Item(109, 12).Value = "Primera"
$ws.Cells.Item(109, 13).Value = 248
$ws.Cells.Item(109, 14).Value = 10000
$ws.Cells.Item(109, 15).Value = 10000
$ws.Cells.Item(109, 16).Value = 10000
$ws.Cells.Item(109, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(109, 18).Value = "Perú"
$ws.Cells.Item(109, 19).Value = 2500
$ws.Cells.Item(109, 20).Value = 4
